$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.8535490096435424
$ws.Range("D2").Value = 0.4025531552197719

$ws.Range("C3").Value = 0.5360848910806454
$ws.Range("D3").Value = 0.597279060083264

$ws.Range("C4").Value = 0.8068510759394214
$ws.Range("D4").Value = 0.4283893044475657

$ws.Range("C5").Value = -0.4917055596603445
$ws.Range("D5").Value = 0.6277977401004358

$ws.Range("C6").Value = 1.125000682153368
$ws.Range("D6").Value = 0.2727133104829305

$ws.Range("C7").Value = 1.810089144151952
$ws.Range("D7").Value = 0.08396434230191518

$ws.Range("C8").Value = 0.5640763809977275
$ws.Range("D8").Value = 0.5784079213336262

$ws.Range("C9").Value = 0.2613969139565138
$ws.Range("D9").Value = 0.7962180175058882

$ws.Range("C10").Value = -0.8250441433187218
$ws.Range("D10").Value = 0.4182031294397579

$ws.Range("C11").Value = -1.119628605206326
$ws.Range("D11").Value = 0.2749473258760742
